$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Naver Blog'
$ws.Range("B2").Value = 20240829
$ws.Range("C2").Value = '강남역고기집 우대포 먹방 유튜버도 반한 소갈비 맛집'
$ws.Range("D2").Value = 'https://blog.naver.com/suzy_blog/223564812454'
$ws.Range("E2").Value = '우대포 강남역점'
$ws.Range("F2").Value = '서울특별시 강남구 강남대로94길 11 1층'
$ws.Range("G2").Value = '#역삼문화공원제1호공영주차장, #우대포, #강남역맛집, #강남역고기집, #역삼동맛집, #강남역회식, #강남역소고기, #강남역고기, #강남역모임, #강남역저녁, #강남역11번출구맛집'
$ws.Range("H2").Value = '2024-08-29 12:52'

# Row 3
$ws.Range("A3").Value = 'Naver Blog'
$ws.Range("B3").Value = 20240829
$ws.Range("C3").Value = '오산 고기 맛집 소오산 오산대역 최고의 소고기 B세트 후기'
$ws.Range("D3").Value = 'https://blog.naver.com/omytheo/223564807189'
$ws.Range("E3").Value = '소오산'
$ws.Range("F3").Value = '경기도 오산시 내삼미로79번길 26 107, 108호'
$ws.Range("G3").Value = '#소오산, #오산고기맛집, #오산대역맛집, #소갈비살, #오산소갈비살, #오산모임, #오산회식, #가성비회식장소, #숯불소고기'
$ws.Range("H3").Value = '2024-08-29 12:52'

# Row 4
$ws.Range("A4").Value = 'Naver Blog'
$ws.Range("B4").Value = 20240829
$ws.Range("C4").Value = '강남N타워 맛집 역삼 파스타바 파브란트'
$ws.Range("D4").Value = 'https://blog.naver.com/winendine2017/223564698837'
$ws.Range("E4").Value = '파브란트'
$ws.Range("F4").Value = '서울특별시 강남구 테헤란로 129 지하2층 B246호'
$ws.Range("G4").Value = '#강남n타워맛집, #강남역파스타, #강남역소개팅, #역삼파스타, #역삼맛집, #파브란트, #강남역데이트, #역삼역데이트, #강남데이트, #강남맛집, #강남파스타, #역삼역파스타, #파스타바, #역삼파스타바, #강남역파스타바'
$ws.Range("H4").Value = '2024-08-29 12:52'

# Row 5
$ws.Range("A5").Value = 'Naver Blog'
$ws.Range("B5").Value = 20240829
$ws.Range("C5").Value = '영덕 강구항 대게 맛집 영덕대게궁 가격'
$ws.Range("D5").Value = 'https://blog.naver.com/ggoulggoul/223564649477'
$ws.Range("E5").Value = '영덕대게궁'
$ws.Range("F5").Value = '경상북도 영덕군 강구면 강구대게길 29 영덕대게궁'
$ws.Range("G5").Value = '#영덕대게맛집, #영덕강구항대게맛집, #강구항대게맛집, #강구항맛집, #영덕대게가격, #영덕대게궁가격, #영덕대게궁, #영덕강구항맛집, #영덕대게, #영덕가볼만한곳, #영덕대게거리, #대게코스요리, #영덕대게거리맛집, #대게宮, #대게궁'
$ws.Range("H5").Value = '2024-08-29 12:52'

# Row 6
$ws.Range("A6").Value = 'Naver Blog'
$ws.Range("B6").Value = 20240829
$ws.Range("C6").Value = '동래 맛집 야키니쿠가 맛있는 현지 감성의 모토이시 공격 후기'
$ws.Range("D6").Value = 'https://blog.naver.com/ganzzajang/223562794470'
$ws.Range("E6").Value = '모토이시 동래점'
$ws.Range("F6").Value = '부산광역시 동래구 충렬대로181번길 65 1~2층'
$ws.Range("G6").Value = '#동래맛집, #동래역맛집, #동래술집, #동래술집추천, #동래고기집, #동래역술집, #명륜동맛집, #동래역고깃집, #부산야끼니꾸, #동래야끼니꾸, #동래와규'
$ws.Range("H6").Value = '2024-08-29 12:52'

# Row 7
$ws.Range("A7").Value = 'Naver Blog'
$ws.Range("B7").Value = 20240829
$ws.Range("C7").Value = '수원 행궁동 한옥 카페 카페도화 연무대점 행리단길 디저트 맛집'
$ws.Range("D7").Value = 'https://blog.naver.com/greenmoom46/223564828431'
$ws.Range("E7").Value = '카페도화 연무대점'
$ws.Range("F7").Value = '경기도 수원시 팔달구 창룡대로103번길 8 2층'
$ws.Range("G7").Value = '#행궁동브런치, #수원브런치, #행궁동와인바, #수원한옥카페, #행궁동한옥카페, #맛띠유, #맛집탐구기'
$ws.Range("H7").Value = '2024-08-29 12:52'

# Row 8
$ws.Range("A8").Value = 'Naver Blog'
$ws.Range("B8").Value = 20240829
$ws.Range("C8").Value = '석촌역 맛집, 청결하고 맛있는 숯불구이 ''2층 통돼지집'' 후기'
$ws.Range("D8").Value = 'https://blog.naver.com/gintonicdays/223564775499'
$ws.Range("E8").Value = '2층 통돼지집'
$ws.Range("F8").Value = '서울특별시 송파구 백제고분로 384 2층 201호(송파동, 평화빌딩)'
$ws.Range("G8").Value = '#2층통돼지집, #잠실맛집, #송파맛집, #석촌역맛집, #돼지두루치기, #돼지숯불구이, #갈비살구이, #테진아제조, #벌집껍데기, #잠실회식장소, #잠실모임장소, #잠실고기맛집, #석촌역고기맛집, #송파동고기집, #송파동고기맛집, #송파동맛집, #롯데타워뷰'
$ws.Range("H8").Value = '2024-08-29 12:52'

# Row 9
$ws.Range("A9").Value = 'Naver Blog'
$ws.Range("B9").Value = 20240829
$ws.Range("C9").Value = '강남역 횟집 속초오징어어시장 가성비 횟집 추천 내돈내산'
$ws.Range("D9").Value = 'https://blog.naver.com/qnwkqortnn/223564770076'
$ws.Range("E9").Value = '속초오징어어시장 강남역점'
$ws.Range("F9").Value = '서울특별시 강남구 강남대로78길 16 1층'
$ws.Range("G9").Value = '#강남역횟집, #강남역가성비횟집, #속초오징어어시장강남역, #강남회식장소추천'
$ws.Range("H9").Value = '2024-08-29 12:52'

# Row 10
$ws.Range("A10").Value = 'Naver Blog'
$ws.Range("B10").Value = 20240829
$ws.Range("C10").Value = '원주 무실동 맛집 추천 점심 물회 부산횟집'
$ws.Range("D10").Value = 'https://blog.naver.com/s2ll27/223564642896'
$ws.Range("E10").Value = '부산횟집'
$ws.Range("F10").Value = '강원특별자치도 원주시 만대공원길 56'
$ws.Range("G10").Value = '#원주무실동맛집, #원주점심맛집, #원주물회맛집, #원주부산횟집, #원주맛집추천, #원주무실동맛집추천'
$ws.Range("H10").Value = '2024-08-29 12:52'

# Row 11
$ws.Range("A11").Value = 'Naver Blog'
$ws.Range("B11").Value = 20240829
$ws.Range("C11").Value = '숙대입구역 맛집 마구로참치 특수부위까지 즐기는 지중해 생참치'
$ws.Range("D11").Value = 'https://blog.naver.com/enviableb/223564632129'
$ws.Range("E11").Value = '마구로참치'
$ws.Range("F11").Value = '서울특별시 용산구 한강대로 271'
$ws.Range("G11").Value = '#숙대입구역맛집, #마구로참치, #숙대입구맛집, #남영역맛집, #갈월동맛집, #용산참치, #남영동맛집, #숙대입구횟집, #남영역횟집, #용산횟집, #용산룸식당, #서울참치맛집'
$ws.Range("H11").Value = '2024-08-29 12:52'

# Row 12
$ws.Range("A12").Value = 'Naver Blog'
$ws.Range("B12").Value = 20240829
$ws.Range("C12").Value = '두툼한 통삼겹이 통째로 김치찌개 의정부맛집 오남매솥밥'
$ws.Range("D12").Value = 'https://blog.naver.com/hahadanro/223562750525'
$ws.Range("E12").Value = '오남매솥밥'
$ws.Range("F12").Value = '경기도 의정부시 오목로225번길 23-3 1층 오남매클라스'
$ws.Range("G12").Value = '#오남매솥밥, #민락동맛집, #한식맛집, #의정부맛집, #의정부김치찌개, #김치찌개, #통삽겹김치찜'
$ws.Range("H12").Value = '2024-08-29 12:52'

# Row 13
$ws.Range("A13").Value = 'Naver Blog'
$ws.Range("B13").Value = 20240829
$ws.Range("C13").Value = '왕림휴게소 근처 동천홍 봉담점 마라 짬뽕 맛집'
$ws.Range("D13").Value = 'https://blog.naver.com/ddomuz/223564273757'
$ws.Range("E13").Value = '동천홍 봉담점'
$ws.Range("F13").Value = '경기도 화성시 봉담읍 삼천병마로 1079-12'
$ws.Range("G13").Value = '#동천홍, #동천홍봉담점, #왕림휴게소맛집, #봉담중국집, #봉담짬뽕맛집, #봉담마라짬뽕, #봉담짜장면, #봉담중식당, #봉담가족모임'
$ws.Range("H13").Value = '2024-08-29 12:52'

# Row 14
$ws.Range("A14").Value = 'Naver Blog'
$ws.Range("B14").Value = 20240829
$ws.Range("C14").Value = '양재역 맛집 미면정 양재본점 만두칼국수 가족모임 추천'
$ws.Range("D14").Value = 'https://blog.naver.com/ykheeee9500/223564513957'
$ws.Range("E14").Value = '미면정 양재본점'
$ws.Range("F14").Value = '서울특별시 강남구 남부순환로365길 22 오이빌딩 2층'
$ws.Range("G14").Value = '#양재역맛집, #양재역맛집미면정, #미면정양재본점, #만두칼국수맛집, #한국음식장인이만드는미면정'
$ws.Range("H14").Value = '2024-08-29 12:52'

# Row 15
$ws.Range("A15").Value = 'Naver Blog'
$ws.Range("B15").Value = 20240829
$ws.Range("C15").Value = '강남구청역 카페 트롬바커피 에그타르트맛집'
$ws.Range("D15").Value = 'https://blog.naver.com/kcy2736/223564554138'
$ws.Range("E15").Value = '트롬바커피'
$ws.Range("F15").Value = '서울특별시 강남구 선릉로126길 12 1층'
$ws.Range("G15").Value = '#서이추환영, #삼성동카페, #강남구청대저트카페, #삼성동디저트카페'
$ws.Range("H15").Value = '2024-08-29 12:52'

# Row 16
$ws.Range("A16").Value = 'Naver Blog'
$ws.Range("B16").Value = 20240829
$ws.Range("C16").Value = '원주시청 근처 근사한 원주 한우 맛집 고기삼촌 후기'
$ws.Range("D16").Value = 'https://blog.naver.com/1115hq/223564798254'
$ws.Range("E16").Value = '고기삼촌'
$ws.Range("F16").Value = '강원특별자치도 원주시 시청로 24 105,106호'
$ws.Range("G16").Value = '#원주한우맛집, #원주시청맛집'
$ws.Range("H16").Value = '2024-08-29 12:52'

# Row 17
$ws.Range("A17").Value = 'Naver Blog'
$ws.Range("B17").Value = 20240829
$ws.Range("C17").Value = '한드케이크_노원 주문제작 예쁜 레터링 케이크 디자인도안 추천'
$ws.Range("D17").Value = 'https://blog.naver.com/sophie242/223561781307'
$ws.Range("E17").Value = '한드케이크'
$ws.Range("F17").Value = '서울특별시 노원구 공릉로58다길 5 B1층 한드케이크'
$ws.Range("G17").Value = '#노원주문제작케이크, #노원레터링케이크, #하계역케이크, #도시락케이크, #한드케이크, #레터링케이크디자인, #주문제작케이크, #주문제작케이크도안'
$ws.Range("H17").Value = '2024-08-29 12:52'

# Row 18
$ws.Range("A18").Value = 'Naver Blog'
$ws.Range("B18").Value = 20240829
$ws.Range("C18").Value = '부천 캠핑식당 이색적인 옥상휴게소'
$ws.Range("D18").Value = 'https://blog.naver.com/buzafama/223561930435'
$ws.Range("E18").Value = '옥상휴게소'
$ws.Range("F18").Value = '경기도 부천시 원미구 부천로 148 1층 전체'
$ws.Range("G18").Value = '#부천캠핑식당, #부천바베큐맛집, #부천캠핑바베큐, #부천캠핑바베큐식당'
$ws.Range("H18").Value = '2024-08-29 12:52'

# Row 19
$ws.Range("A19").Value = 'Naver Blog'
$ws.Range("B19").Value = 20240829
$ws.Range("C19").Value = '대구 수성구 술집 고퀄리티 안주에 분위기 좋은 범어동 맛집 한주'
$ws.Range("D19").Value = 'https://blog.naver.com/ejr8806/223564731432'
$ws.Range("E19").Value = '한주'
$ws.Range("F19").Value = '대구광역시 수성구 달구벌대로512길 10 1층 한주'
$ws.Range("G19").Value = '#대구수성구술집, #수성구술집, #한주, #범어동한주, #범어동맛집'
$ws.Range("H19").Value = '2024-08-29 12:52'

# Row 20
$ws.Range("A20").Value = 'Naver Blog'
$ws.Range("B20").Value = 20240829
$ws.Range("C20").Value = '판교 술집 하이네켄 생맥주 맛집 챠오바라이트'
$ws.Range("D20").Value = 'https://blog.naver.com/darkyoon85/223564335342'
$ws.Range("E20").Value = '챠오바라이트 판교점'
$ws.Range("F20").Value = '경기도 성남시 분당구 대왕판교로 660 유스페이스1 1층 125호 챠오바라이트'
$ws.Range("G20").Value = '#판교맛집, #판교술집, #하이네켄, #하이네켄생맥주, #하이네켄생맥주맛집, #챠오바라이트, #생맥주, #첫모금, #맛집추천, #광고'
$ws.Range("H20").Value = '2024-08-29 12:52'

# Row 21
$ws.Range("A21").Value = 'Naver Blog'
$ws.Range("B21").Value = 20240829
$ws.Range("C21").Value = '연남동 국숯집 남발카세 팝업 - 젊고 내공 뛰어난 남발게 김준승 셰프님의 새롭고 기발하고 훌륭한 해산물 요리 잔치'
$ws.Range("D21").Value = 'https://blog.naver.com/zephyr122059/223559250507'
$ws.Range("E21").Value = '국숯집 연남점'
$ws.Range("F21").Value = '서울특별시 마포구 동교로34길 21 우측 지층'
$ws.Range("G21").Value = '#남발게님, #국숯집, #연남동국숯집, #남발카세, #해산물잔치, #국숯집팝업, #남발게님팝업, #연남동남발카세팝업'
$ws.Range("H21").Value = '2024-08-29 12:52'
